# Applies the Aug 28 2023 crypto price/volume refresh described by the commit diff.
# Row 12/13 (WrappedEther <-> Polkadot) swapped rank position, so B/C/D/E are all rewritten there.
# Cells whose new text looks like a plain number (e.g. "219.37") are entered with a leading
# apostrophe -- exactly as a user would type them in Excel -- so they stay text instead of
# being auto-converted to a General-formatted number (which would also eat trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.218.19'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '1.655.71'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('E4').Value = '  -0.66%  '
$ws.Range('D5').Value = '''219.37'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').Value = '''0.5245'
$ws.Range('E6').Value = '  -1.72%  '
$ws.Range('D7').Value = '''1.004'
$ws.Range('E7').Value = '  -0.69%  '
$ws.Range('D8').Value = '''0.2674'
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').Value = '''0.06386'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('D10').Value = '''20.63'
$ws.Range('E10').Value = '  -1.56%  '
$ws.Range('D11').Value = '''0.07744'
$ws.Range('E11').Value = '  -1.49%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '''4.597'
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.676.90'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('D14').Value = '1.883.34'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('D15').Value = '''0.5658'
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('D16').Value = '0.0₅8318'
$ws.Range('E16').Value = '  +1.55%  '
$ws.Range('D17').Value = '''65.53'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '26.212.79'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').Value = '''1.005'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D20').Value = '''4.706'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').Value = '''192.39'
$ws.Range('E21').Value = '  -3.50%  '
$ws.Range('D22').Value = '''10.40'
$ws.Range('E22').Value = '  +0.57%  '
$ws.Range('D23').Value = '''6.038'
$ws.Range('E23').Value = '  -0.51%  '
$ws.Range('D24').Value = '''1.005'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').Value = '''143.51'
$ws.Range('E25').Value = '  -2.21%  '
$ws.Range('D26').Value = '''0.1201'
$ws.Range('D27').Value = '''7.287'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('D28').Value = '''15.99'
$ws.Range('E28').Value = '  -1.84%  '
$ws.Range('D29').Value = '''1.500'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').Value = '''0.05631'
$ws.Range('D31').Value = '''1.277'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').Value = '''3.501'
$ws.Range('E32').Value = '  -1.67%  '
$ws.Range('D33').Value = '''3.381'
$ws.Range('E33').Value = '  +2.12%  '
$ws.Range('D34').Value = '''1.583'
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('D35').Value = '''2.803'
$ws.Range('E35').Value = '  -1.39%  '
$ws.Range('D36').Value = '''0.9473'
$ws.Range('E36').Value = '  -2.30%  '
$ws.Range('D37').Value = '''2.406'
$ws.Range('E37').Value = '  -1.36%  '
$ws.Range('D38').Value = '''0.5760'
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('D39').Value = '''0.01624'
$ws.Range('E39').Value = '  +0.64%  '
$ws.Range('D40').Value = '''5.909'
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('D41').Value = '''2.575'
$ws.Range('E41').Value = '  -0.42%  '
$ws.Range('D42').Value = '''0.8463'
$ws.Range('E42').Value = '  -2.04%  '
$ws.Range('E43').Value = '  -0.71%  '
$ws.Range('D44').Value = '1.028.56'
$ws.Range('E44').Value = '  -4.87%  '
$ws.Range('D45').Value = '''101.76'
$ws.Range('E45').Value = '  -2.12%  '
$ws.Range('D46').Value = '1.793.57'
$ws.Range('E46').Value = '  -1.00%  '
$ws.Range('D47').Value = '''58.52'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '0.0₈105'
$ws.Range('E48').Value = '  +3.03%  '
$ws.Range('D49').Value = '''1.004'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('D50').Value = '''0.05326'
$ws.Range('E50').Value = '  +3.07%  '
$ws.Range('D51').Value = '''8.029'
$ws.Range('E51').Value = '  +0.36%  '
